# "penambahan metode dan pembaruan pada penelitian terdahulu"
#
# 1. Insert a new "Sheet1" worksheet right after "Progress 2" (before "Note"),
#    listing the related-field categories.
# 2. On "Progress 2", fill in the Metode (K) / Keterkaitan Penelitian (L) /
#    Pembaruan (M) columns for rows 3-8, and mark row 3's Status (N) as DONE.
# 3. Restore the original active sheet / selections when finished.

$wb = $excel.ActiveWorkbook

$progress2 = $wb.Worksheets.Item("Progress 2")

# --- 1. New "Sheet1" worksheet, inserted right after "Progress 2" ---------
$sheet1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $progress2)

# Write L3 first so the shared "Sama sama face detection..." string becomes
# the very next new shared-string entry before the Sheet1 literals.
$progress2.Range("L3").Value = "Sama sama face detection untuk drowsiness"

$sheet1.Range("C3").Value = "RPL"
$sheet1.Range("C4").Value = "UI UX"
$sheet1.Range("C5").Value = "BIOMEDICAL"
$sheet1.Range("C6").Value = "DATA MINING"
$sheet1.Range("C8").Value = "MULTI MEDIA PROCESSING"
$sheet1.Range("C9").Value = "Komunikasi Visual"
$sheet1.Columns.Item(3).ColumnWidth = 24
$sheet1.Range("B8:D9").Select()

# --- 2. "Progress 2" Metode / Pembaruan columns ----------------------------
$progress2.Range("K3").Value = "pelatihan algoritma Haar Cascade untuk pengenalan wajah. Algoritma ini kemudian diterapkan pada aplikasi berbasis Android untuk mendeteksi wajah pengemudi secara otomatis. Sistem ini bertujuan untuk memantau tanda-tanda kelelahan pada pengemudi dalam konteks logistik dan transportasi umum, sehingga dapat meningkatkan keselamatan jalan raya."
$progress2.Range("M3").Value = "embaruan dalam penelitian ini meliputi penerapan sistem menggunakan aplikasi Android yang terintegrasi dengan server, memungkinkan pemantauan kondisi pengemudi secara real-time. Integrasi ini bertujuan untuk meningkatkan efektivitas deteksi kelelahan pengemudi pada aplikasi berbasis Android, mendukung pengawasan dan peningkatan keselamatan dalam lalu lantas"

$progress2.Range("K4").Value = " penggunaan algoritma YOLOv3 untuk deteksi objek. Algoritma ini kemudian digabungkan dengan LSTM (Long Short-Term Memory) dalam proses pelatihan data untuk menganalisis perilaku temporal pengemudi. Kombinasi ini bertujuan untuk meningkatkan akurasi deteksi tanda-tanda kelelahan pada pengemudi."
$progress2.Range("M4").Value = "peningkatan efektivitas sistem melalui integrasi algoritma YOLOv3 dengan LSTM. LSTM digunakan untuk menganalisis data urutan waktu (temporal data) dari perilaku pengemudi, memungkinkan sistem untuk mendeteksi tanda-tanda kelelahan dengan lebih akurat berdasarkan pola perilaku yang berubah seiring waktu."

$progress2.Range("M5").Value = "penerapan YOLOv5 yang dibantu dengan Vision Transformers untuk meningkatkan akurasi klasifikasi. Integrasi Vision Transformers membantu sistem dalam memahami dan menganalisis fitur visual secara lebih mendalam, sehingga meningkatkan kemampuan deteksi kantuk pada pengemudi."
$progress2.Range("K5").Value = "penerapan YOLOv5 untuk deteksi objek. Algoritma YOLOv5 kemudian digabungkan dengan Vision Transformers untuk menganalisis dan mengklasifikasikan fitur visual dari wajah dan perilaku pengemudi. Kombinasi ini bertujuan untuk meningkatkan akurasi deteksi tanda-tanda kelelahan pada pengemudi"

$progress2.Range("M6").Value = "penggunaan website serta bot Telegram untuk membantu mengirimkan informasi kepada pengemudi ataupun pengawas. Integrasi ini memungkinkan pemantauan dan pemberian peringatan secara real-time, sehingga meningkatkan respons terhadap tanda-tanda kelelahan pada pengemudi."
$progress2.Range("K6").Value = "pelatihan YOLOv4 dengan dataset RAKIBUL.ECE.RUET. Algoritma YOLOv4 diterapkan untuk mendeteksi ekspresi wajah yang menunjukkan tanda-tanda kelelahan pada pengemudi secara real-time melalui aplikasi berbasis web. Kombinasi ini bertujuan untuk meningkatkan akurasi dan efisiensi deteksi kelelahan pada pengemudi."

$progress2.Range("M7").Value = "melakukan uji coba variasi pelatihan untuk memperoleh akurasi maksimal. Dengan menguji berbagai konfigurasi dan parameter pelatihan, penelitian ini bertujuan untuk mengoptimalkan kinerja algoritma YOLOv4 dalam mendeteksi tanda-tanda kantuk pada pengemudi."
$progress2.Range("K7").Value = "pelatihan algoritma YOLOv4 menggunakan dataset open source. Algoritma YOLOv4 diterapkan untuk mendeteksi ekspresi wajah yang menunjukkan tanda-tanda kantuk pada pengemudi. Penelitian ini fokus pada pengujian berbagai variasi pelatihan untuk meningkatkan akurasi deteksi kantuk secara keseluruhan."

$progress2.Range("M8").Value = "optimalisasi pelatihan menggunakan CUDA dan Adam optimizer untuk memaksimalkan efisiensi pelatihan. Penggunaan CUDA memungkinkan pemrosesan paralel pada GPU, sementara Adam optimizer meningkatkan kecepatan konvergensi dan stabilitas model selama pelatihan."
$progress2.Range("K8").Value = "pelatihan YOLOv5 agar dapat melakukan deteksi wajah. Algoritma YOLOv5 diterapkan untuk mendeteksi ekspresi wajah yang menunjukkan tanda-tanda kantuk pada pengemudi. Penelitian ini mengoptimalkan proses pelatihan dengan menggunakan CUDA dan Adam optimizer untuk meningkatkan efisiensi dan akurasi deteksi."

# Same "Keterkaitan Penelitian" note for the remaining rows (re-uses the
# shared string created above for L3).
$progress2.Range("L4").Value = "Sama sama face detection untuk drowsiness"
$progress2.Range("L5").Value = "Sama sama face detection untuk drowsiness"
$progress2.Range("L6").Value = "Sama sama face detection untuk drowsiness"
$progress2.Range("L7").Value = "Sama sama face detection untuk drowsiness"
$progress2.Range("L8").Value = "Sama sama face detection untuk drowsiness"

# Row 3 is now finished/reviewed.
$progress2.Range("N3").Value = "DONE"

# --- 3. View bookkeeping ----------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("F40").Select()

# Leave "Progress 2" as the active sheet/selection, matching the original.
$progress2.Activate()
$progress2.Range("M13").Select()
